$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3434516666666667
$ws.Range("H2").Value = 1.030355
$ws.Range("I2").Value = 0.07181106008319683
$ws.Range("J2").Value = 0.07181106008319683
$ws.Range("M2").Value = 2.819352333333333
$ws.Range("N2").Value = 8.458057
$ws.Range("O2").Value = 0.05609595354453337
$ws.Range("P2").Value = 0.05609595354453337
$ws.Range("Q2").Value = 0.9683112578038888
$ws.Range("R2").Value = 8.714801320234999
$ws.Range("S2").Value = 0.004028309890410704
$ws.Range("T2").Value = 0.004028309890410704
$ws.Range("G3").Value = 0.3434516666666667
$ws.Range("H3").Value = 1.030355
$ws.Range("I3").Value = 0.07181106008319683
$ws.Range("J3").Value = 0.07181106008319683
$ws.Range("O3").Value = 0.7103421657667366
$ws.Range("P3").Value = 0.7103421657667366
$ws.Range("Q3").Value = 12.26171002617278
$ws.Range("R3").Value = 110.355390235555
$ws.Range("S3").Value = 0.05101042394550329
$ws.Range("T3").Value = 0.05101042394550329
$ws.Range("G4").Value = 0.3434516666666667
$ws.Range("H4").Value = 1.030355
$ws.Range("I4").Value = 0.07181106008319683
$ws.Range("J4").Value = 0.07181106008319683
$ws.Range("M4").Value = 11.738694
$ws.Range("N4").Value = 35.216082
$ws.Range("O4").Value = 0.23356188068873
$ws.Range("P4").Value = 0.23356188068873
$ws.Range("Q4").Value = 4.03167401879
$ws.Range("R4").Value = 36.28506616911
$ws.Range("S4").Value = 0.01677232624728284
$ws.Range("T4").Value = 0.01677232624728284
$ws.Range("I5").Value = 0.3084341675553999
$ws.Range("J5").Value = 0.3084341675553999
$ws.Range("M5").Value = 2.819352333333333
$ws.Range("N5").Value = 8.458057
$ws.Range("O5").Value = 0.05609595354453337
$ws.Range("P5").Value = 0.05609595354453337
$ws.Range("Q5").Value = 4.158973233221334
$ws.Range("R5").Value = 37.43075909899201
$ws.Range("S5").Value = 0.01730190873473453
$ws.Range("T5").Value = 0.01730190873473453
$ws.Range("I6").Value = 0.3084341675553999
$ws.Range("J6").Value = 0.3084341675553999
$ws.Range("O6").Value = 0.7103421657667366
$ws.Range("P6").Value = 0.7103421657667366
$ws.Range("S6").Value = 0.2190937945777633
$ws.Range("T6").Value = 0.2190937945777633
$ws.Range("I7").Value = 0.3084341675553999
$ws.Range("J7").Value = 0.3084341675553999
$ws.Range("M7").Value = 11.738694
$ws.Range("N7").Value = 35.216082
$ws.Range("O7").Value = 0.23356188068873
$ws.Range("P7").Value = 0.23356188068873
$ws.Range("Q7").Value = 17.316357931488
$ws.Range("R7").Value = 155.847221383392
$ws.Range("S7").Value = 0.07203846424290207
$ws.Range("T7").Value = 0.07203846424290207
$ws.Range("G8").Value = 2.964109000000001
$ws.Range("H8").Value = 8.892327000000002
$ws.Range("I8").Value = 0.6197547723614032
$ws.Range("J8").Value = 0.6197547723614033
$ws.Range("M8").Value = 2.819352333333333
$ws.Range("N8").Value = 8.458057
$ws.Range("O8").Value = 0.05609595354453337
$ws.Range("P8").Value = 0.05609595354453337
$ws.Range("Q8").Value = 8.356867625404334
$ws.Range("R8").Value = 75.21180862863902
$ws.Range("S8").Value = 0.03476573491938813
$ws.Range("T8").Value = 0.03476573491938813
$ws.Range("G9").Value = 2.964109000000001
$ws.Range("H9").Value = 8.892327000000002
$ws.Range("I9").Value = 0.6197547723614032
$ws.Range("J9").Value = 0.6197547723614033
$ws.Range("O9").Value = 0.7103421657667366
$ws.Range("P9").Value = 0.7103421657667366
$ws.Range("Q9").Value = 105.8228815620897
$ws.Range("R9").Value = 952.4059340588072
$ws.Range("S9").Value = 0.44023794724347
$ws.Range("T9").Value = 0.44023794724347
$ws.Range("G10").Value = 2.964109000000001
$ws.Range("H10").Value = 8.892327000000002
$ws.Range("I10").Value = 0.6197547723614032
$ws.Range("J10").Value = 0.6197547723614033
$ws.Range("M10").Value = 11.738694
$ws.Range("N10").Value = 35.216082
$ws.Range("O10").Value = 0.23356188068873
$ws.Range("P10").Value = 0.23356188068873
$ws.Range("Q10").Value = 34.79476853364601
$ws.Range("R10").Value = 313.152916802814
$ws.Range("S10").Value = 0.1447510901985451
$ws.Range("T10").Value = 0.1447510901985451
